# EnergyFlowsMath.xlsx update
# - Rename the "e_mat_Recycled_HQ" label (cell A22 on Sheet1) to
#   "e_mat_Recycled_HQ + LQ" to reflect that recycled HQ and LQ energies
#   are now combined in the energy-flow calculation.
# - Leave the active cell selection on A23, matching where the author's
#   cursor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A22").Value = "e_mat_Recycled_HQ + LQ"

$ws.Range("A23").Select()
